$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the DuckDB/SQL queries stored in column B (rows 2-7) and in C2.
# The join conditions were changed from the generic ".id" columns to the
# renamed, table-specific id columns (e.g. "study_id", "participant_id").
# ---------------------------------------------------------------------------
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $text = [string]$ws.Range($addr).Value()
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $ws.Range($addr).Value = $text
}

# ---------------------------------------------------------------------------
# Widen column C so the (now longer) query text in C2 fits without relying
# on "best fit" auto-sizing.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = (200 / 3)
